# SG2042EVB config.xlsx - "power" sheet: add CHECK/PG_* port rows (checkport for 2042)
#
# The table of enable/power-good signals grew from 16 data rows (A2:D17) to
# 20 data rows (A2:D21): several "PG_*" (CHECK) rows were inserted and the
# remaining rows shifted down, picking up the net names/types/delays that
# used to sit one row below them. Net effect: rewrite A2:D21 in full.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("power")

# Name | Net name | Type | Delay
$data = @(
    @("EN_VDD_3V3",         "EN_VDD_3V3",         "ENABLE", "1000"),
    @("EN_VDDIO18",         "EN_VDDIO18",         "ENABLE", "1000"),
    @("EN_VDDC",            "EN_VDDC",            "ENABLE", "1000"),
    @("PG_VDDC",            "PG_VDDC",            "CHECK",  "0"),
    @("EN_VDDIO33",         "EN_VDDIO33",         "ENABLE", "1000"),
    @("EN_DDR_VDD_0V8",     "EN_DDR_VDD_0V8",     "ENABLE", "0"),
    @("EN_VDD_PCIE_D_0V8",  "EN_VDD_PCIE_D_0V8",  "ENABLE", "0"),
    @("EN_VDD_PLL_0V8",     "EN_VDD_PLL_0V8",     "ENABLE", "0"),
    @("PG_DDR_VDD_0V8",     "PG_DDR_VDD_0V8",     "CHECK",  "0"),
    @("EN_VDD_PCIE_H_1V8",  "EN_VDD_PCIE_H_1V8",  "ENABLE", "1000"),
    @("EN_DDR01_VPP_2V5",   "EN_DDR01_VPP_2V5",   "ENABLE", "0"),
    @("EN_DDR23_VPP_2V5",   "EN_DDR23_VPP_2V5",   "ENABLE", "1000"),
    @("EN_DDR01_VDDQ_1V2",  "EN_DDR01_VDDQ_1V2",  "ENABLE", "0"),
    @("EN_DDR23_VDDQ_1V2",  "EN_DDR23_VDDQ_1V2",  "ENABLE", "1000"),
    @("PG_DDR01_VDDQ_1V2",  "PG_DDR01_VDDQ_1V2",  "CHECK",  "0"),
    @("PG_DDR23_VDDQ_1V2",  "PG_DDR23_VDDQ_1V2",  "CHECK",  "0"),
    @("EN_DDR01_VTT_0V6",   "EN_DDR01_VTT_0V6",   "ENABLE", "0"),
    @("EN_DDR23_VTT_0V6",   "EN_DDR23_VTT_0V6",   "ENABLE", "1000"),
    @("EN_VQPS18",          "EN_VQPS18",          "ENABLE", "30000"),
    @("SYS_RST_X_H",        "SYS_RST_X_H",        "ENABLE", "1000")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]

    # "Delay" keeps being stored as text (matches the source workbook, where
    # the whole A1:E column range carries numberStoredAsText) rather than a
    # genuine number, so force Text formatting across the round-trip without
    # leaving a lingering custom cell style behind.
    $dcell = $ws.Cells.Item($row, 4)
    $dcell.NumberFormat = "@"
    $dcell.Value = $data[$i][3]
    $dcell.ClearFormats()
}
